$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Adjust column widths for columns B, C and D
    $ws.Columns.Item(2).ColumnWidth = 18.8333333333333
    $ws.Columns.Item(3).ColumnWidth = 14.8333333333333
    $ws.Columns.Item(4).ColumnWidth = 30.8333333333333

    # Bump the format version shown in C2
    $ws.Range("C2").Value = "v0.1.3"
}
